$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the column headers: "<Name>_old" -> "<Name>_FV2404", "<Name>_new" -> "<Name>_FV2410" ---
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $leftCol  = $i + 1     # columns A..J
    $rightCol = $i + 12    # columns L..U
    $ws.Cells.Item(1, $leftCol).Value  = $baseNames[$i] + "_FV2404"
    $ws.Cells.Item(1, $rightCol).Value = $baseNames[$i] + "_FV2410"
}

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

# --- 3. Turn the used range into an Excel Table ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$lo.Name = "Table1"
$ws.Range("A1").Select() | Out-Null

Write-Host "done"
